$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "CD36"
$ws.Range("B4").Value = "long-chain fatty acid and lipoprotein presception"
$ws.Range("C4").Value = "https://academic.oup.com/hmg/article/17/11/1695/599714"
